$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "SQ"
$ws.Range("A4").Value = "AAPL"

$ws.Range("A4").Select()
